$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h) hold plain text values such as
# "22.149.30", "0.9998" or "  -0.74%  " rather than numbers/dates. A plain
# Range.Value assignment would let Excel auto-convert numeric-looking text
# (e.g. "1.000" -> 1, "0.9998" -> 0.9998 as a Number) so each new value is
# written with a leading apostrophe, Excel's standard "force text" entry,
# which keeps the string exactly as scraped.

$ws.Range('D2').Value = "'22.149.30"
$ws.Range('E2').Value = "'  -0.74%  "
$ws.Range('D3').Value = "'1.558.75"
$ws.Range('E3').Value = "'  +0.06%  "
$ws.Range('D4').Value = "'0.9983"
$ws.Range('E4').Value = "'  -0.24%  "
$ws.Range('D5').Value = "'0.9998"
$ws.Range('E5').Value = "'  -0.09%  "
$ws.Range('D6').Value = "'291.82"
$ws.Range('E6').Value = "'  +1.51%  "
$ws.Range('D7').Value = "'0.3945"
$ws.Range('E7').Value = "'  +4.39%  "
$ws.Range('D8').Value = "'0.3244"
$ws.Range('E8').Value = "'  -0.72%  "
$ws.Range('D9').Value = "'43.77"
$ws.Range('E9').Value = "'  -0.10%  "
$ws.Range('D10').Value = "'0.07335"
$ws.Range('E10').Value = "'  -0.45%  "
$ws.Range('D11').Value = "'1.091"
$ws.Range('E11').Value = "'  -4.09%  "
$ws.Range('D12').Value = "'0.9982"
$ws.Range('E12').Value = "'  -0.25%  "
$ws.Range('D13').Value = "'19.15"
$ws.Range('E13').Value = "'  -5.46%  "
$ws.Range('D14').Value = "'5.671"
$ws.Range('D15').Value = "'0.00001145"
$ws.Range('E15').Value = "'  +6.52%  "
$ws.Range('D16').Value = "'6.672"
$ws.Range('E16').Value = "'  -1.24%  "
$ws.Range('D17').Value = "'1.555.99"
$ws.Range('E17').Value = "'  -0.94%  "
$ws.Range('D18').Value = "'0.06610"
$ws.Range('E18').Value = "'  -0.43%  "
$ws.Range('D19').Value = "'83.91"
$ws.Range('E19').Value = "'  -2.44%  "
$ws.Range('D20').Value = "'1.000"
$ws.Range('E20').Value = "'  -0.12%  "
$ws.Range('D21').Value = "'6.326"
$ws.Range('E21').Value = "'  -1.03%  "
$ws.Range('E22').Value = "'  -1.55%  "
$ws.Range('D23').Value = "'11.33"
$ws.Range('E23').Value = "'  -2.62%  "
$ws.Range('D24').Value = "'22.153.39"
$ws.Range('E24').Value = "'  -0.73%  "
$ws.Range('D25').Value = "'2.338"
$ws.Range('E25').Value = "'  +2.11%  "
$ws.Range('D26').Value = "'2.447"
$ws.Range('E26').Value = "'  -4.20%  "
$ws.Range('D27').Value = "'148.20"
$ws.Range('E27').Value = "'  -1.86%  "
$ws.Range('D28').Value = "'18.69"
$ws.Range('E28').Value = "'  -3.15%  "
$ws.Range('D29').Value = "'4.877"
$ws.Range('E29').Value = "'  -1.07%  "
$ws.Range('D30').Value = "'1.728.31"
$ws.Range('E30').Value = "'  -0.96%  "
$ws.Range('D31').Value = "'119.21"
$ws.Range('E31').Value = "'  -2.66%  "
$ws.Range('D32').Value = "'1.029"
$ws.Range('E32').Value = "'  -4.41%  "
$ws.Range('D33').Value = "'5.737"
$ws.Range('E33').Value = "'  -2.67%  "
$ws.Range('D34').Value = "'0.08367"
$ws.Range('E34').Value = "'  +1.72%  "
$ws.Range('D35').Value = "'1.624"
$ws.Range('E35').Value = "'  -13.82%  "
$ws.Range('D36').Value = "'9.115"
$ws.Range('E36').Value = "'  -2.92%  "
$ws.Range('D37').Value = "'0.06164"
$ws.Range('E37').Value = "'  -2.00%  "
$ws.Range('D38').Value = "'0.02281"
$ws.Range('E38').Value = "'  -3.72%  "
$ws.Range('D39').Value = "'5.164"
$ws.Range('E39').Value = "'  -2.39%  "
$ws.Range('D40').Value = "'1.214"
$ws.Range('E40').Value = "'  -1.99%  "
$ws.Range('D41').Value = "'0.2068"
$ws.Range('E41').Value = "'  -3.77%  "
$ws.Range('D42').Value = "'10.83"
$ws.Range('E42').Value = "'  -1.96%  "
$ws.Range('D43').Value = "'0.9997"
$ws.Range('E43').Value = "'  -0.15%  "
$ws.Range('D44').Value = "'0.5889"
$ws.Range('E44').Value = "'  -2.84%  "
$ws.Range('D45').Value = "'13.10"
$ws.Range('E45').Value = "'  -4.36%  "
$ws.Range('D46').Value = "'3.766"
$ws.Range('E46').Value = "'  +0.66%  "
$ws.Range('D47').Value = "'0.5648"
$ws.Range('E47').Value = "'  -4.55%  "
$ws.Range('D48').Value = "'118.84"
$ws.Range('E48').Value = "'  -3.63%  "
$ws.Range('D49').Value = "'1.904"
$ws.Range('E49').Value = "'  -3.74%  "
$ws.Range('D50').Value = "'1.143"
$ws.Range('E50').Value = "'  -2.87%  "
$ws.Range('D51').Value = "'0.06876"
$ws.Range('E51').Value = "'  -2.88%  "
